$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 2 (shifts existing data rows 2-21 down to 8-27)
$ws.Rows("2:7").Insert()
$ws.Rows("2:7").ClearFormats()

$newTopRows = @(
    @(0.0357356183230876, 0.0120645882561802, 0.1313360333442688),
    @(0.0088575463742017, 0.0383317954838275, 0.0606283769011497),
    @(-0.0204639863222837, 0.0259617734700441, 0.0542142912745475),
    @(-0.0181732401251792, 0.0203112699091434, -0.0135917514562606),
    @(-0.0335975885391235, -0.0102319931611418, -0.0829249545931816),
    @(-0.0200058370828628, -0.0244346093386411, -0.0164933614432811)
)

$r = 2
foreach ($row in $newTopRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$newBottomRows = @(
    @(0.0577267669141292, 0.2557998299598694, 0.0554360225796699),
    @(-0.0160352122038602, -0.030695978552103, -0.0510072484612464),
    @(-0.0320704244077205, -0.107512280344963, -0.04505131021142),
    @(0.0154243474826216, 0.1117883399128913, 0.0210748501121997)
)

$r = 28
foreach ($row in $newBottomRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
